$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value as TEXT (prevents Excel from auto-coercing
# numeric-looking strings, e.g. "512.41" or "1.00", into numbers/floats).
function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" "69.590.82"
$ws.Range("E2").Value = "  +1.88%  "

Set-TextValue "D3" "3.957.81"
$ws.Range("E3").Value = "  +0.81%  "

$ws.Range("E4").Value = "  +0.03%  "

Set-TextValue "D5" "512.41"
$ws.Range("E5").Value = "  +5.54%  "

Set-TextValue "D6" "147.70"
$ws.Range("E6").Value = "  +0.23%  "

Set-TextValue "D7" "0.627"
$ws.Range("E7").Value = "  +0.39%  "

$ws.Range("E8").Value = "  +0.03%  "

Set-TextValue "D9" "0.737"
$ws.Range("E9").Value = "  +0.55%  "

$ws.Range("E10").Value = "  +5.33%  "

Set-TextValue "D11" "0.0000349"
$ws.Range("E11").Value = "  -0.54%  "

Set-TextValue "D12" "43.63"
$ws.Range("E12").Value = "  +1.66%  "

Set-TextValue "D13" "10.53"
$ws.Range("E13").Value = "  -1.62%  "

Set-TextValue "D14" "4.580.07"
$ws.Range("E14").Value = "  +0.71%  "

Set-TextValue "D15" "3.964.97"
$ws.Range("E15").Value = "  +0.89%  "

$ws.Range("E16").Value = "  -1.19%  "

$ws.Range("E17").Value = "  -0.28%  "

$ws.Range("E18").Value = "  +7.63%  "

$ws.Range("E19").Value = "  +0.53%  "

Set-TextValue "D20" "69.663.65"
$ws.Range("E20").Value = "  +1.82%  "

Set-TextValue "D21" "437.24"
$ws.Range("E21").Value = "  -1.18%  "

Set-TextValue "D22" "3.45"
$ws.Range("E22").Value = "  -1.17%  "

Set-TextValue "D23" "14.70"
$ws.Range("E23").Value = "  -2.56%  "

Set-TextValue "D24" "89.05"
$ws.Range("E24").Value = "  +0.80%  "

Set-TextValue "D25" "11.87"
$ws.Range("E25").Value = "  +5.78%  "

Set-TextValue "D26" "3.89"
$ws.Range("E26").Value = "  +7.59%  "

Set-TextValue "D27" "11.24"
$ws.Range("E27").Value = "  -3.11%  "

Set-TextValue "D28" "37.28"
$ws.Range("E28").Value = "  -4.47%  "

Set-TextValue "D29" "5.67"
$ws.Range("E29").Value = "  -2.08%  "

Set-TextValue "D30" "709.30"
$ws.Range("E30").Value = "  -1.21%  "

Set-TextValue "D31" "13.43"
$ws.Range("E31").Value = "  -2.47%  "

$ws.Range("E32").Value = "  -1.04%  "

$ws.Range("E33").Value = "  -0.60%  "

Set-TextValue "D34" "66.19"
$ws.Range("E34").Value = "  +8.33%  "

Set-TextValue "D35" "0.447"
$ws.Range("E35").Value = "  +11.46%  "

Set-TextValue "D36" "0.0₃0884"

Set-TextValue "D37" "6.07"
$ws.Range("E37").Value = "  -5.00%  "

Set-TextValue "D38" "40.96"
$ws.Range("E38").Value = "  -3.29%  "

$ws.Range("E39").Value = "  +1.18%  "

Set-TextValue "D40" "1.00"
$ws.Range("E40").Value = "  +0.02%  "

$ws.Range("E41").Value = "  -0.13%  "

$ws.Range("E42").Value = "  +1.91%  "

Set-TextValue "D43" "2.91"
$ws.Range("E43").Value = "  -1.70%  "

Set-TextValue "D44" "3.12"
$ws.Range("E44").Value = "  +6.66%  "

$ws.Range("E45").Value = "  -4.52%  "

Set-TextValue "D46" "0.145"
$ws.Range("E46").Value = "  +1.81%  "

$ws.Range("E47").Value = "  +3.15%  "

$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue "D48" "0.0₆0359"
$ws.Range("E48").Value = "  +1.13%  "

$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D49" "3.01"
$ws.Range("E49").Value = "  +6.21%  "

$ws.Range("B50").Value = "LidoDAOToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue "D50" "3.41"
$ws.Range("E50").Value = "  -0.23%  "

$ws.Range("E51").Value = "  -0.95%  "
